# This script applies the edit described by the commit "Fruta / hortaliza, semanal".
# A new weekly price record is inserted as row 85 of the sheet, pushing the
# previously existing rows 85-148 down to rows 86-149 (and therefore also
# appending a brand-new last row, 149, that is simply the old row 148's data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 85; Excel automatically shifts rows
# 85..148 down to 86..149, carrying their values and formatting with them.
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new weekly record.
$ws.Cells.Item(85, 1).Value = 10
$ws.Cells.Item(85, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(85, 3).Value = "La Araucanía"
$ws.Cells.Item(85, 4).Value = 44824
$ws.Cells.Item(85, 5).Value = 9
$ws.Cells.Item(85, 6).Value = "Fruta"
$ws.Cells.Item(85, 7).Value = 100107
$ws.Cells.Item(85, 8).Value = "Otros"
$ws.Cells.Item(85, 9).Value = 100107002
$ws.Cells.Item(85, 10).Value = "Chirimoya"
$ws.Cells.Item(85, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(85, 12).Value = "Primera"
$ws.Cells.Item(85, 13).Value = 55
$ws.Cells.Item(85, 14).Value = 3500
$ws.Cells.Item(85, 15).Value = 3600
$ws.Cells.Item(85, 16).Value = 3536
$ws.Cells.Item(85, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(85, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(85, 19).Value = 3536
$ws.Cells.Item(85, 20).Value = 1
